$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# 1) Insert a new "setup" worksheet immediately before "heats"
# ------------------------------------------------------------------
$heatsSheet = $wb.Worksheets.Item("heats")
$setupSheet = $wb.Worksheets.Add($heatsSheet)
$setupSheet.Name = "setup"

# Populate the new "setup" sheet
$setupSheet.Range("A1").Value = "Calorimeter"
$setupSheet.Range("B1").Value = "ampoule"
$setupSheet.Range("A2").Value = "Initial volume"
$setupSheet.Range("B2").Value = 1
$setupSheet.Range("A3").Value = "component"
$setupSheet.Range("B3").Value = "L"

# ------------------------------------------------------------------
# 2) Rework the "targets" sheet: drop the old layout, keep only the
#    "constants " label in A1
# ------------------------------------------------------------------
$targetsSheet = $wb.Worksheets.Item("targets")
$targetsSheet.Rows("2:2").Delete()
$targetsSheet.Cells.ClearContents()
$targetsSheet.Range("A1").Value = "constants "

# ------------------------------------------------------------------
# 3) Make "setup" the active sheet / active selection (A4)
# ------------------------------------------------------------------
$setupSheet.Activate()
$setupSheet.Range("A4").Select()
